# edit.ps1
# Applies the benchmark-results update + README note + algorithm execution
# sequence change described in the commit message / xml diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Updated raw benchmark numbers in columns C (row backtracking),
#    D (block backtracking) and E (deduction) for puzzles in rows 8-18.
#    (Only the cells whose numbers actually changed are touched; unchanged
#    values - e.g. D9, C10, D10, E15:E18 - are left alone.)
# ---------------------------------------------------------------------------
$ws.Range("C8").Value  = 7380
$ws.Range("D8").Value  = 8819
$ws.Range("E8").Value  = 11198

$ws.Range("C9").Value  = 382383
$ws.Range("E9").Value  = 17909

$ws.Range("E10").Value = 37858

$ws.Range("C11").Value = 632639
$ws.Range("D11").Value = 798050
$ws.Range("E11").Value = 106140

$ws.Range("C12").Value = 275644
$ws.Range("D12").Value = 306227
$ws.Range("E12").Value = 81753

$ws.Range("C13").Value = 360593
$ws.Range("D13").Value = 216975
$ws.Range("E13").Value = 138219

$ws.Range("C14").Value = 5796629
$ws.Range("D14").Value = 9114752
$ws.Range("E14").Value = 80804

$ws.Range("C15").Value = 86150620
$ws.Range("D15").Value = 47439333

$ws.Range("C16").Value = 3336116
$ws.Range("D16").Value = 1257266

$ws.Range("C17").Value = 2097809
$ws.Range("D17").Value = 1611088

$ws.Range("C18").Value = 7955443
$ws.Range("D18").Value = 19118997

# ---------------------------------------------------------------------------
# 2. Algorithm execution sequence changed: column H ("puzzle number" on the
#    right-hand comparison table) used to hold static numbers (and, for row
#    14, a wrong hard-coded 11). Now each cell mirrors column B with a
#    formula, fixing that bug along the way.
# ---------------------------------------------------------------------------
$ws.Range("H8").Formula  = "=B8"
$ws.Range("H9:H18").Formula = "=B9"

# ---------------------------------------------------------------------------
# 3. The "ALL LEVELS" summary row no longer shows a computed average for the
#    deduction-algorithm column (it is marked "x" like the per-level rows
#    that didn't complete), and the derived speed-up ratios L25/M25 are
#    cleared out instead of being computed.
# ---------------------------------------------------------------------------
$ws.Range("K24").Copy()
$ws.Range("K25").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("K25").Value = "x"

$ws.Range("L25").ClearContents()
$ws.Range("M25").ClearContents()

# ---------------------------------------------------------------------------
# 4. README note: the "Conclusion:" heading and the final concluding
#    sentence are removed from the notes section below the tables.
# ---------------------------------------------------------------------------
$ws.Range("A27").ClearContents()
$ws.Range("A28").ClearContents()

# ---------------------------------------------------------------------------
# 5. Force a full recalculation so every dependent formula (I/J/K columns,
#    the averages in rows 21-25, and the speed-up ratios) reflects the new
#    raw data.
# ---------------------------------------------------------------------------
$excel.CalculateFull()

# ---------------------------------------------------------------------------
# 6. Restore the scroll position / selection that was active when the
#    author last saved the workbook.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3   # "C1" -> column C is left-most visible
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K19").Select()
